$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "92.691.70"
$ws.Range("E2").Value = "  +6.14%  "

$ws.Range("D3").Value = "3.296.75"
$ws.Range("E3").Value = "  +1.19%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "636.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.83%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.416"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +14.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.727"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.16%  "

$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("D10").Value = "3.294.08"
$ws.Range("E10").Value = "  +1.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.594"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.87%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.180"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.88%  "

$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000264"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.59%  "

$ws.Range("D15").Value = "3.891.15"
$ws.Range("E15").Value = "  +0.91%  "

$ws.Range("D16").Value = "92.159.61"
$ws.Range("E16").Value = "  +5.80%  "

$ws.Range("E17").Value = "  +1.46%  "

$ws.Range("D18").Value = "3.294.07"
$ws.Range("E18").Value = "  +1.10%  "

$ws.Range("E19").Value = "  +6.60%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.45%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "440.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.90%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000189"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +47.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.22%  "

$ws.Range("D27").Value = "3.476.40"
$ws.Range("E27").Value = "  +1.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "76.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.50%  "

$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.179"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "559.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.67"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +25.68%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.73"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.72%  "

$ws.Range("E39").Value = "  -3.12%  "

$ws.Range("E40").Value = "  +3.66%  "

$ws.Range("E41").Value = "  +0.10%  "

$ws.Range("E42").Value = "  +1.59%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.06%  "

$ws.Range("E44").Value = "  +0.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "151.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "181.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.70%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.131"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.47%  "

$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "44.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.01%  "

$ws.Range("E49").Value = "  +0.78%  "

$ws.Range("E50").Value = "  +2.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.90%  "
